$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("M2").Value = 0.02648366666666667
$ws.Range("N2").Value = 0.07945099999999999
$ws.Range("O2").Value = 0.001430039273477916
$ws.Range("P2").Value = 0.001430039273477917
$ws.Range("Q2").Value = 0.0006235755874444445
$ws.Range("R2").Value = 0.005612180287
$ws.Range("S2").Value = 0.000004026354787907393
$ws.Range("T2").Value = 0.000004026354787907393

# Row 3
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("O3").Value = 0.7016741634339546
$ws.Range("P3").Value = 0.7016741634339547
$ws.Range("Q3").Value = 0.3059684351142222
$ws.Range("R3").Value = 2.753715916028
$ws.Range("S3").Value = 0.001975602474624517
$ws.Range("T3").Value = 0.001975602474624517

# Row 4
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.2968957972925674
$ws.Range("P4").Value = 0.2968957972925675
$ws.Range("Q4").Value = 0.129462857867
$ws.Range("R4").Value = 1.165165720803
$ws.Range("S4").Value = 0.0008359265630734945
$ws.Range("T4").Value = 0.0008359265630734945

# Row 5
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("M5").Value = 0.02648366666666667
$ws.Range("N5").Value = 0.07945099999999999
$ws.Range("O5").Value = 0.001430039273477916
$ws.Range("P5").Value = 0.001430039273477917
$ws.Range("Q5").Value = 0.2185618000394444
$ws.Range("R5").Value = 1.967056200355
$ws.Range("S5").Value = 0.001411228033555558
$ws.Range("T5").Value = 0.001411228033555558

# Row 6
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("O6").Value = 0.7016741634339546
$ws.Range("P6").Value = 0.7016741634339547
$ws.Range("S6").Value = 0.6924440945257248
$ws.Range("T6").Value = 0.6924440945257248

# Row 7
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.2968957972925674
$ws.Range("P7").Value = 0.2968957972925675
$ws.Range("S7").Value = 0.2929903254790364
$ws.Range("T7").Value = 0.2929903254790365

# Row 8
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("M8").Value = 0.02648366666666667
$ws.Range("N8").Value = 0.07945099999999999
$ws.Range("O8").Value = 0.001430039273477916
$ws.Range("P8").Value = 0.001430039273477917
$ws.Range("Q8").Value = 0.002289786647888889
$ws.Range("R8").Value = 0.020608079831
$ws.Range("S8").Value = 0.00001478488513445089
$ws.Range("T8").Value = 0.00001478488513445089

# Row 9
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("O9").Value = 0.7016741634339546
$ws.Range("P9").Value = 0.7016741634339547
$ws.Range("S9").Value = 0.007254466433605358
$ws.Range("T9").Value = 0.007254466433605358

# Row 10
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.2968957972925674
$ws.Range("P10").Value = 0.2968957972925675
$ws.Range("S10").Value = 0.003069545250457494
$ws.Range("T10").Value = 0.003069545250457495
